$d = $word.ActiveDocument

# 1. Replace "כל רעיון בקובץ נפרד" with "בגיטהאב של האירגון שלכם" inside the
#    parenthetical remark in the weekly-homework paragraph, e.g.
#    "(כל רעיון בקובץ נפרד) והוסיפו " -> "(בגיטהאב של האירגון שלכם) והוסיפו "
$found = $d.Content.Find.Execute("כל רעיון בקובץ נפרד", $false, $false, $false, $false, $false, $true, 1, $false, "בגיטהאב של האירגון שלכם", 2)

# 2. Remove the leftover "_GoBack" bookmark (Word adds this automatically to
#    mark the last edit location; it is stripped on a clean save).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$d.Saved = $false
